$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:D51 as text first so numeric-looking price strings
# (e.g. "0.994", "569.82") are stored as text, matching the source data,
# not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.000.33'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '2.443.34'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '569.82'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = '146.39'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '2.446.05'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').Value = '0.111'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '0.356'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '5.22'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '26.77'
$ws.Range('E14').Value = '  +2.46%  '
$ws.Range('D15').Value = '0.0000179'
$ws.Range('E15').Value = '  +3.22%  '
$ws.Range('D16').Value = '2.876.14'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '63.128.18'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = '2.452.00'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = '11.32'
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').Value = '7.19'
$ws.Range('E20').Value = '  +5.47%  '
$ws.Range('D21').Value = '323.49'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = '4.16'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '0.994'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '1.92'
$ws.Range('E24').Value = '  +10.94%  '
$ws.Range('D25').Value = '66.40'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').Value = '615.34'
$ws.Range('E26').Value = '  +10.46%  '
$ws.Range('D27').Value = '8.64'
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').Value = '0.0000102'
$ws.Range('E28').Value = '  +9.64%  '
$ws.Range('D29').Value = '2.590.66'
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').Value = '1.47'
$ws.Range('E31').Value = '  +5.72%  '
$ws.Range('D32').Value = '8.22'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Value = '0.143'
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').Value = '5.08'
$ws.Range('E35').Value = '  +7.46%  '
$ws.Range('D36').Value = '1.50'
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '0.381'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = '18.62'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '5.35'
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '145.49'
$ws.Range('E41').Value = '  -4.54%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('D43').Value = '2.60'
$ws.Range('E43').Value = '  +14.69%  '
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('D45').Value = '146.85'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').Value = '3.71'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').Value = '0.0538'
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('D48').Value = '20.57'
$ws.Range('E48').Value = '  +3.42%  '
$ws.Range('D49').Value = '0.600'
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('D50').Value = '0.0233'
$ws.Range('E50').Value = '  +2.27%  '
$ws.Range('D51').Value = '0.0919'
$ws.Range('E51').Value = '  -0.17%  '

# Restore the default (General) formatting on the price column so no
# stray per-cell style lingers now that the text values are locked in.
$ws.Range("D2:D51").ClearFormats()

Write-Output "Applied cryptos update."
